$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENTRADAS")

$ws.Range("A1").Value = "S112"
$ws.Range("B1").Value = "S110"
$ws.Range("C1").Value = "S109"
